$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.122774
$ws.Range("H2").Value = 33.368322
$ws.Range("I2").Value = 0.2449652610853511
$ws.Range("J2").Value = 0.2449652610853511
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 30.99161333333333
$ws.Range("N2").Value = 92.97484
$ws.Range("O2").Value = 0.3599121977633812
$ws.Range("P2").Value = 0.3599121977633811
$ws.Range("Q2").Value = 344.7127110020533
$ws.Range("R2").Value = 3102.41439901848
$ws.Range("S2").Value = 0.0881659854929092
$ws.Range("T2").Value = 0.08816598549290917

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.122774
$ws.Range("H3").Value = 33.368322
$ws.Range("I3").Value = 0.2449652610853511
$ws.Range("J3").Value = 0.2449652610853511
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 29.913269
$ws.Range("N3").Value = 89.739807
$ws.Range("O3").Value = 0.3473891556493311
$ws.Range("P3").Value = 0.3473891556493311
$ws.Range("Q3").Value = 332.718530688206
$ws.Range("R3").Value = 2994.466776193854
$ws.Range("S3").Value = 0.08509827521185807
$ws.Range("T3").Value = 0.08509827521185806

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.122774
$ws.Range("H4").Value = 33.368322
$ws.Range("I4").Value = 0.2449652610853511
$ws.Range("J4").Value = 0.2449652610853511
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 19.150218
$ws.Range("N4").Value = 57.450654
$ws.Range("O4").Value = 0.2223955550134164
$ws.Range("P4").Value = 0.2223955550134163
$ws.Range("Q4").Value = 213.003546864732
$ws.Range("R4").Value = 1917.031921782588
$ws.Range("S4").Value = 0.05447918519808311
$ws.Range("T4").Value = 0.0544791851980831

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.122774
$ws.Range("H5").Value = 33.368322
$ws.Range("I5").Value = 0.2449652610853511
$ws.Range("J5").Value = 0.2449652610853511
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.053716000000001
$ws.Range("N5").Value = 18.161148
$ws.Range("O5").Value = 0.07030309157387134
$ws.Range("P5").Value = 0.07030309157387132
$ws.Range("Q5").Value = 67.33411492818401
$ws.Range("R5").Value = 606.007034353656
$ws.Range("S5").Value = 0.01722181518250074
$ws.Range("T5").Value = 0.01722181518250073

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 18.220714
$ws.Range("H6").Value = 54.662142
$ws.Range("I6").Value = 0.4012885600454987
$ws.Range("J6").Value = 0.4012885600454988
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 30.99161333333333
$ws.Range("N6").Value = 92.97484
$ws.Range("O6").Value = 0.3599121977633812
$ws.Range("P6").Value = 0.3599121977633811
$ws.Range("Q6").Value = 564.6893229452534
$ws.Range("R6").Value = 5082.20390650728
$ws.Range("S6").Value = 0.144428647583278
$ws.Range("T6").Value = 0.144428647583278

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 18.220714
$ws.Range("H7").Value = 54.662142
$ws.Range("I7").Value = 0.4012885600454987
$ws.Range("J7").Value = 0.4012885600454988
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 29.913269
$ws.Range("N7").Value = 89.739807
$ws.Range("O7").Value = 0.3473891556493311
$ws.Range("P7").Value = 0.3473891556493311
$ws.Range("Q7").Value = 545.0411192540661
$ws.Range("R7").Value = 4905.370073286595
$ws.Range("S7").Value = 0.1394032940459417
$ws.Range("T7").Value = 0.1394032940459417

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 18.220714
$ws.Range("H8").Value = 54.662142
$ws.Range("I8").Value = 0.4012885600454987
$ws.Range("J8").Value = 0.4012885600454988
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 19.150218
$ws.Range("N8").Value = 57.450654
$ws.Range("O8").Value = 0.2223955550134164
$ws.Range("P8").Value = 0.2223955550134163
$ws.Range("Q8").Value = 348.930645215652
$ws.Range("R8").Value = 3140.375806940868
$ws.Range("S8").Value = 0.08924479203185337
$ws.Range("T8").Value = 0.08924479203185336

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 18.220714
$ws.Range("H9").Value = 54.662142
$ws.Range("I9").Value = 0.4012885600454987
$ws.Range("J9").Value = 0.4012885600454988
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.053716000000001
$ws.Range("N9").Value = 18.161148
$ws.Range("O9").Value = 0.07030309157387134
$ws.Range("P9").Value = 0.07030309157387132
$ws.Range("Q9").Value = 110.303027873224
$ws.Range("R9").Value = 992.7272508590161
$ws.Range("S9").Value = 0.02821182638442566
$ws.Range("T9").Value = 0.02821182638442566

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1189986666666667
$ws.Range("H10").Value = 0.356996
$ws.Range("I10").Value = 0.002620797603979787
$ws.Range("J10").Value = 0.002620797603979787
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 30.99161333333333
$ws.Range("N10").Value = 92.97484
$ws.Range("O10").Value = 0.3599121977633812
$ws.Range("P10").Value = 0.3599121977633811
$ws.Range("Q10").Value = 3.687960664515555
$ws.Range("R10").Value = 33.19164598064
$ws.Range("S10").Value = 0.0009432570255413687
$ws.Range("T10").Value = 0.0009432570255413684

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1189986666666667
$ws.Range("H11").Value = 0.356996
$ws.Range("I11").Value = 0.002620797603979787
$ws.Range("J11").Value = 0.002620797603979787
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 29.913269
$ws.Range("N11").Value = 89.739807
$ws.Range("O11").Value = 0.3473891556493311
$ws.Range("P11").Value = 0.3473891556493311
$ws.Range("Q11").Value = 3.559639126641333
$ws.Range("R11").Value = 32.036752139772
$ws.Range("S11").Value = 0.0009104366667743282
$ws.Range("T11").Value = 0.0009104366667743281

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1189986666666667
$ws.Range("H12").Value = 0.356996
$ws.Range("I12").Value = 0.002620797603979787
$ws.Range("J12").Value = 0.002620797603979787
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 19.150218
$ws.Range("N12").Value = 57.450654
$ws.Range("O12").Value = 0.2223955550134164
$ws.Range("P12").Value = 0.2223955550134163
$ws.Range("Q12").Value = 2.278850408376
$ws.Range("R12").Value = 20.509653675384
$ws.Range("S12").Value = 0.0005828537377149165
$ws.Range("T12").Value = 0.0005828537377149163

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1189986666666667
$ws.Range("H13").Value = 0.356996
$ws.Range("I13").Value = 0.002620797603979787
$ws.Range("J13").Value = 0.002620797603979787
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 6.053716000000001
$ws.Range("N13").Value = 18.161148
$ws.Range("O13").Value = 0.07030309157387134
$ws.Range("P13").Value = 0.07030309157387132
$ws.Range("Q13").Value = 0.7203841323786666
$ws.Range("R13").Value = 6.483457191408
$ws.Range("S13").Value = 0.0001842501739491735
$ws.Range("T13").Value = 0.0001842501739491735

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 15.943029
$ws.Range("H14").Value = 47.829087
$ws.Range("I14").Value = 0.3511253812651704
$ws.Range("J14").Value = 0.3511253812651704
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 30.99161333333333
$ws.Range("N14").Value = 92.97484
$ws.Range("O14").Value = 0.3599121977633812
$ws.Range("P14").Value = 0.3599121977633811
$ws.Range("Q14").Value = 494.10019013012
$ws.Range("R14").Value = 4446.90171117108
$ws.Range("S14").Value = 0.1263743076616526
$ws.Range("T14").Value = 0.1263743076616526

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 15.943029
$ws.Range("H15").Value = 47.829087
$ws.Range("I15").Value = 0.3511253812651704
$ws.Range("J15").Value = 0.3511253812651704
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 29.913269
$ws.Range("N15").Value = 89.739807
$ws.Range("O15").Value = 0.3473891556493311
$ws.Range("P15").Value = 0.3473891556493311
$ws.Range("Q15").Value = 476.908115151801
$ws.Range("R15").Value = 4292.173036366209
$ws.Range("S15").Value = 0.121977149724757
$ws.Range("T15").Value = 0.121977149724757

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 15.943029
$ws.Range("H16").Value = 47.829087
$ws.Range("I16").Value = 0.3511253812651704
$ws.Range("J16").Value = 0.3511253812651704
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 19.150218
$ws.Range("N16").Value = 57.450654
$ws.Range("O16").Value = 0.2223955550134164
$ws.Range("P16").Value = 0.2223955550134163
$ws.Range("Q16").Value = 305.312480930322
$ws.Range("R16").Value = 2747.812328372898
$ws.Range("S16").Value = 0.07808872404576502
$ws.Range("T16").Value = 0.07808872404576499

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 15.943029
$ws.Range("H17").Value = 47.829087
$ws.Range("I17").Value = 0.3511253812651704
$ws.Range("J17").Value = 0.3511253812651704
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 6.053716000000001
$ws.Range("N17").Value = 18.161148
$ws.Range("O17").Value = 0.07030309157387134
$ws.Range("P17").Value = 0.07030309157387132
$ws.Range("Q17").Value = 96.51456974576402
$ws.Range("R17").Value = 868.631127711876
$ws.Range("S17").Value = 0.02468519983299576
$ws.Range("T17").Value = 0.02468519983299576

